# Weekly refresh of the Camote (Vega Modelo de Temuco) price records.
# The 21 existing daily price rows are re-shuffled into a new date order and
# one brand-new record is appended (row 23), matching the upstream weekly
# export. Columns A,B,C,E,F,G,H,I,R are constant across every record in this
# sub-sheet, so only D (Fecha) and J..Q (Volumen..Kg o Unidades) vary.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$constA = 10
$constB = 'Vega Modelo de Temuco'
$constC = 'La Araucanía'
$constE = 9
$constF = 100114002
$constG = 'Camote'
$constH = 'Sin especificar'
$constI = 'Primera'
$constR = 'Hortaliza'

$rows = @(
    @{Row=2; D=44424; J=30; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=3; D=44188; J=20; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=4; D=44364; J=15; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Perú'; P=1333; Q=15},
    @{Row=5; D=44389; J=45; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=6; D=44425; J=10; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=7; D=44385; J=18; K=20000; L=20000; M=20000; N='$/malla 20 kilos'; O='Región de Arica y Parinacota'; P=1000; Q=20},
    @{Row=8; D=44340; J=40; K=18000; L=18000; M=18000; N='$/malla 20 kilos'; O='Perú'; P=900; Q=20},
    @{Row=9; D=44186; J=20; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=10; D=44179; J=20; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=11; D=44321; J=15; K=25000; L=25000; M=25000; N='$/caja 15 kilos granel'; O='Perú'; P=1667; Q=15},
    @{Row=12; D=44316; J=20; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=13; D=44341; J=40; K=17000; L=18000; M=17500; N='$/malla 20 kilos'; O='Perú'; P=875; Q=20},
    @{Row=14; D=44441; J=40; K=20000; L=20000; M=20000; N='$/malla 20 kilos'; O='Perú'; P=1000; Q=20},
    @{Row=15; D=44369; J=20; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=16; D=44369; J=20; K=20000; L=20000; M=20000; N='$/malla 20 kilos'; O='Región de Arica y Parinacota'; P=1000; Q=20},
    @{Row=17; D=44294; J=5; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Perú'; P=1333; Q=15},
    @{Row=18; D=44315; J=30; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=19; D=44315; J=30; K=20000; L=20000; M=20000; N='$/malla 20 kilos'; O='Región de Arica y Parinacota'; P=1000; Q=20},
    @{Row=20; D=44329; J=40; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Perú'; P=1333; Q=15},
    @{Row=21; D=44161; J=20; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=22; D=44438; J=40; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15},
    @{Row=23; D=44175; J=20; K=20000; L=20000; M=20000; N='$/caja 15 kilos granel'; O='Región de Arica y Parinacota'; P=1333; Q=15}
)

# Column D (Fecha) carries a custom date-time number format on the existing
# rows (style index 2 in the original workbook); remember it so the new row
# we append (row 23) gets the same formatting instead of the default style.
$dateFormat = $ws.Range("D2").NumberFormat

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $constI
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
    $ws.Cells.Item($r, 12).Value = $item.L
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 14).Value = $item.N
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
    $ws.Cells.Item($r, 17).Value = $item.Q
    $ws.Cells.Item($r, 18).Value = $constR

    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
}
